$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.340.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.89%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.940.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.91%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.07%  "

# Row 6
$ws.Range("E6").Value = "  -5.93%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3349"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07324"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.87%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.75%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08132"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.939.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.551"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.352.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008304"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "255.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.892"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.194.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.966"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.859"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.424"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.20%  "

# Row 29
$ws.Range("E29").Value = "  -9.35%  "

# Row 30
$ws.Range("E30").Value = "  -3.51%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.345"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.474"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.258"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05243"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.275"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7596"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.742"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.44%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02009"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "

# Row 40
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.681"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4575"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.06%  "

# Row 43
$ws.Range("E43").Value = "  -5.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8437"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.61%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.909"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.442"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.92%  "

# Row 49
$ws.Range("E49").Value = "  +0.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("E51").Value = "  -4.00%  "

Write-Host "Update complete"